$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-10-10 12:34:56"

# Update the "取得日時" (fetched-at) timestamp for existing rows 2-12.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# New row 13
$ws.Range("A13").Value = $newTimestamp
$ws.Range("B13").Value = "【急募】salamに関するウェブサイト制作の依頼"
$ws.Range("C13").Value = "システム開発"
$ws.Range("D13").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E13").Value = "期限情報なし"
$ws.Range("F13").Value = "https://www.lancers.jp/work/detail/5411046"
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5411046")
$ws.Range("G13").Value = 38
$ws.Range("H13").Value = "◇サイト"

# New row 14
$ws.Range("A14").Value = $newTimestamp
$ws.Range("B14").Value = "【急募】16タイプ診断コンテンツのLP制作"
$ws.Range("C14").Value = "システム開発"
$ws.Range("D14").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E14").Value = "期限情報なし"
$ws.Range("F14").Value = "https://www.lancers.jp/work/detail/5408735"
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5408735")
$ws.Range("G14").Value = 25

# New row 15
$ws.Range("A15").Value = $newTimestamp
$ws.Range("B15").Value = "〖リモート可〗Delphiエンジニア募集"
$ws.Range("C15").Value = "システム開発"
$ws.Range("D15").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E15").Value = "期限情報なし"
$ws.Range("F15").Value = "https://www.lancers.jp/work/detail/5341051"
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.lancers.jp/work/detail/5341051")
$ws.Range("G15").Value = 25

# New row 16
$ws.Range("A16").Value = $newTimestamp
$ws.Range("B16").Value = "初回 【継続あり】Microsoft PL-300/400/600 資格試験向け問題集作成"
$ws.Range("C16").Value = "システム開発"
$ws.Range("D16").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E16").Value = "期限情報なし"
$ws.Range("F16").Value = "https://www.lancers.jp/work/detail/5411149"
$ws.Hyperlinks.Add($ws.Range("F16"), "https://www.lancers.jp/work/detail/5411149")
$ws.Range("G16").Value = 18

# New row 17
$ws.Range("A17").Value = $newTimestamp
$ws.Range("B17").Value = "急募 【緊急】selenium(ruby)でのX自動ログインの実装"
$ws.Range("C17").Value = "システム開発"
$ws.Range("D17").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E17").Value = "期限情報なし"
$ws.Range("F17").Value = "https://www.lancers.jp/work/detail/5411088"
$ws.Hyperlinks.Add($ws.Range("F17"), "https://www.lancers.jp/work/detail/5411088")
$ws.Range("G17").Value = 18

# Apply the Hyperlink style to the new URL cells only after every hyperlink
# has been created, so the engine reuses the existing "Hyperlink" cell style
# (same as F2:F12) instead of registering a duplicate style record.
$ws.Range("F13").Style = "Hyperlink"
$ws.Range("F14").Style = "Hyperlink"
$ws.Range("F15").Style = "Hyperlink"
$ws.Range("F16").Style = "Hyperlink"
$ws.Range("F17").Style = "Hyperlink"
